$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 additions (A2 already has value "(SEO ou accessiblité ?)")
$ws.Range("B2").Value = "JS non minifier"
$ws.Range("C2").Value = "analyse rapide index.html"

# Row 3
$ws.Range("A3").Value = "seo"
$ws.Range("B3").Value = "meta description & title"

# Row 4
$ws.Range("A4").Value = "seo"
$ws.Range("B4").Value = "meta keywords ?"

# Row 5
$ws.Range("A5").Value = "seo"
$ws.Range("B5").Value = "balise title ??"

# Row 6
$ws.Range("A6").Value = "seo"
$ws.Range("B6").Value = "!!! Black Hat mot nav bar : class keywords "

# Row 7
$ws.Range("A7").Value = "??"
$ws.Range("B7").Value = "liens 404 errors"

# Row 8
$ws.Range("A8").Value = "access"
$ws.Range("B8").Value = "langue balise"

# Row 9
$ws.Range("A9").Value = "access"
$ws.Range("B9").Value = "revoir le responsive"

# Update the active selection to B10, matching the diff
$ws.Range("B10").Select()
